# This script applies the daily cryptos-list price/volume refresh described
# by the commit "Updated cryptos list ... with GitHub Actions".
#
# Columns:
#   B = Coin name, C = Link, D = Price, E = Volume(1h)
#
# Some Price values are plain numeric-looking strings (e.g. "6.60", "4.10").
# Assigning such a string straight to Range.Value causes Excel to parse it as
# a number and normalize it (losing the trailing zero / exact text), so for
# those specific values we prefix the string with a leading apostrophe to
# force a literal-text interpretation, exactly preserving the source text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.581.54'
$ws.Range('E2').Value = '  +0.31%  '
$ws.Range('D3').Value = '3.120.89'
$ws.Range('E3').Value = '  -1.42%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = '''571.74'
$ws.Range('E5').Value = '  +0.07%  '
$ws.Range('D6').Value = '''162.95'
$ws.Range('E6').Value = '  -3.40%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('E8').Value = '  -5.15%  '
$ws.Range('D9').Value = '3.124.15'
$ws.Range('E9').Value = '  -1.91%  '
$ws.Range('E10').Value = '  -1.52%  '
$ws.Range('D11').Value = '''6.60'
$ws.Range('E11').Value = '  -3.16%  '
$ws.Range('E12').Value = '  -2.15%  '
$ws.Range('D13').Value = '3.662.40'
$ws.Range('E13').Value = '  -1.42%  '
$ws.Range('E14').Value = '  -2.31%  '
$ws.Range('D15').Value = '64.617.29'
$ws.Range('E15').Value = '  +0.23%  '
$ws.Range('D16').Value = '''24.79'
$ws.Range('E16').Value = '  -2.24%  '
$ws.Range('D17').Value = '3.131.62'
$ws.Range('E17').Value = '  -1.80%  '
$ws.Range('E18').Value = '  -1.27%  '
$ws.Range('D19').Value = '''408.41'
$ws.Range('E19').Value = '  -2.15%  '
$ws.Range('D20').Value = '''5.22'
$ws.Range('E20').Value = '  -1.71%  '
$ws.Range('D21').Value = '''12.39'
$ws.Range('D22').Value = '''6.99'
$ws.Range('E22').Value = '  -2.17%  '
$ws.Range('D23').Value = '''0.999'
$ws.Range('E23').Value = '  -0.14%  '
$ws.Range('D24').Value = '''68.01'
$ws.Range('E24').Value = '  -2.44%  '
$ws.Range('D25').Value = '''0.479'
$ws.Range('E25').Value = '  -3.85%  '
$ws.Range('E26').Value = '  -5.15%  '
$ws.Range('D27').Value = '''0.0000103'
$ws.Range('E27').Value = '  -1.01%  '
$ws.Range('D28').Value = '''9.13'
$ws.Range('E28').Value = '  +4.21%  '
$ws.Range('E29').Value = '  -0.21%  '
$ws.Range('D30').Value = '''0.999'
$ws.Range('E30').Value = '  +0.07%  '
$ws.Range('E31').Value = '  -1.34%  '
$ws.Range('E32').Value = '  -2.38%  '
$ws.Range('D33').Value = '''164.23'
$ws.Range('E33').Value = '  +5.58%  '
$ws.Range('E34').Value = '  -2.59%  '
$ws.Range('B35').Value = 'Aptos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D35').Value = '''6.22'
$ws.Range('E35').Value = '  -2.19%  '
$ws.Range('B36').Value = 'Fetch.AI'
$ws.Range('C36').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D36').Value = '''1.13'
$ws.Range('E36').Value = '  +0.63%  '
$ws.Range('E37').Value = '  -0.77%  '
$ws.Range('E38').Value = '  -2.14%  '
$ws.Range('D39').Value = '2.597.58'
$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D40').Value = '''4.10'
$ws.Range('E40').Value = '  -2.60%  '
$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D41').Value = '''23.62'
$ws.Range('E41').Value = '  -1.89%  '
$ws.Range('D42').Value = '''38.29'
$ws.Range('E42').Value = '  -1.97%  '
$ws.Range('D43').Value = '''0.688'
$ws.Range('E43').Value = '  -4.09%  '
$ws.Range('E44').Value = '  -0.56%  '
$ws.Range('D45').Value = '''5.23'
$ws.Range('E45').Value = '  -3.99%  '
$ws.Range('D46').Value = '''0.0255'
$ws.Range('E46').Value = '  -3.32%  '
$ws.Range('B47').Value = 'Bittensor'
$ws.Range('C47').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D47').Value = '''285.94'
$ws.Range('E47').Value = '  -1.37%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').Value = '''21.05'
$ws.Range('E48').Value = '  -1.82%  '
$ws.Range('D49').Value = '''0.995'
$ws.Range('E49').Value = '  -0.51%  '
$ws.Range('D50').Value = '''0.0972'
$ws.Range('E50').Value = '  -1.91%  '
$ws.Range('E51').Value = '  +0.32%  '
